# Update the "users" sheet with data extracted/parsed via Google Sheets.
# Row 2 is replaced with a fresh "NewUser" account, and two more student
# rows are appended (rows 3 and 4), replacing the old admin/blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray empty D3 cell left over from the previous layout so the
# sheet's used range goes back down to A1:C4.
$ws.Range("D3").ClearContents()

# Row 2: NewUser / Qwertz / 20ФиПЛ-1
$ws.Range("A2").Value = "NewUser"
$ws.Range("B2").Value = "Qwertz"
$ws.Range("C2").Value = "20ФиПЛ-1"

# Row 3: Pervokursnik / 1 / 23ФИЛ-1
$ws.Range("A3").Value = "Pervokursnik"
$ws.Range("B3").Value = "1"
$ws.Range("C3").Value = "23ФИЛ-1"

# Row 4: Student / Student / 22ФиПЛ-1
$ws.Range("A4").Value = "Student"
$ws.Range("B4").Value = "Student"
$ws.Range("C4").Value = "22ФиПЛ-1"

# Move the active selection like the author's last save did.
$ws.Range("D8").Select()
